$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "tool_ID == “custody_answer %}", $true, $false, $false, $false, $false,
    $true, 1, $false, "tool_ID == “custody_answer” %}", 2)

$d.Content.Find.Execute(
    "tool_ID == “divorce_answer %}", $true, $false, $false, $false, $false,
    $true, 1, $false, "tool_ID == “divorce_answer” %}", 2)
